$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "A1393201"
$ws.Range("D2").Value = "SMatrix_Salm_5523"
$ws.Range("E2").Value = "20220304-TestAut-PA-35523"
$ws.Range("R2").Value = "'03/04/2022"
$ws.Range("T2").Value = "TestCartridge5523"

# Row 3
$ws.Range("A3").Value = "A1393202"
$ws.Range("D3").Value = "SMatrix_Salm_5523"
$ws.Range("E3").Value = "20220304-TestAut-PA-35523"
$ws.Range("R3").Value = "'03/04/2022"
$ws.Range("T3").Value = "TestCartridge5523"

# Row 4
$ws.Range("A4").Value = "A1393203"
$ws.Range("D4").Value = "SMatrix_Salm_5523"
$ws.Range("E4").Value = "20220304-TestAut-PA-35523"
$ws.Range("R4").Value = "'03/04/2022"
$ws.Range("T4").Value = "TestCartridge5523"

# Row 5
$ws.Range("A5").Value = "A1393204"
$ws.Range("D5").Value = "SMatrix_Salm_5523"
$ws.Range("E5").Value = "20220304-TestAut-PA-35523"
$ws.Range("R5").Value = "'03/04/2022"
$ws.Range("T5").Value = "TestCartridge5523"

# Row 6
$ws.Range("A6").Value = "A1393205"
$ws.Range("D6").Value = "SMatrix_Salm_5523"
$ws.Range("E6").Value = "20220304-TestAut-PA-35523"
$ws.Range("R6").Value = "'03/04/2022"
$ws.Range("T6").Value = "TestCartridge5523"

# Row 7
$ws.Range("A7").Value = "A1393209"
$ws.Range("D7").Value = "SMatrix_Salm_5523"
$ws.Range("E7").Value = "20220304-TestAut-PA-35523"
$ws.Range("Q7").Value = "'9"
$ws.Range("R7").Value = "'03/04/2022"
$ws.Range("T7").Value = "TestCartridge5523"

# Row 8
$ws.Range("A8").Value = "A1393210"
$ws.Range("D8").Value = "SMatrix_Salm_5523"
$ws.Range("E8").Value = "20220304-TestAut-PA-35523"
$ws.Range("Q8").Value = "'10"
$ws.Range("R8").Value = "'03/04/2022"
$ws.Range("T8").Value = "TestCartridge5523"

# Row 9
$ws.Range("A9").Value = "A1393211"
$ws.Range("D9").Value = "SMatrix_Salm_5523"
$ws.Range("E9").Value = "20220304-TestAut-PA-35523"
$ws.Range("Q9").Value = "'11"
$ws.Range("R9").Value = "'03/04/2022"
$ws.Range("T9").Value = "TestCartridge5523"

# Row 10
$ws.Range("A10").Value = "A1393206"
$ws.Range("D10").Value = "SMatrix_Salm_5523"
$ws.Range("E10").Value = "20220304-TestAut-PA-35523"
$ws.Range("Q10").Value = "'6"
$ws.Range("R10").Value = "'03/04/2022"
$ws.Range("T10").Value = "TestCartridge5523"

# Row 11
$ws.Range("A11").Value = "A1393207"
$ws.Range("D11").Value = "SMatrix_Salm_5523"
$ws.Range("E11").Value = "20220304-TestAut-PA-35523"
$ws.Range("Q11").Value = "'7"
$ws.Range("R11").Value = "'03/04/2022"
$ws.Range("T11").Value = "TestCartridge5523"

# Row 12
$ws.Range("A12").Value = "A1393208"
$ws.Range("D12").Value = "SMatrix_Salm_5523"
$ws.Range("E12").Value = "20220304-TestAut-PA-35523"
$ws.Range("Q12").Value = "'8"
$ws.Range("R12").Value = "'03/04/2022"
$ws.Range("T12").Value = "TestCartridge5523"

# Row 13
$ws.Range("A13").Value = "A1393212"
$ws.Range("D13").Value = "SMatrix_Salm_5523"
$ws.Range("E13").Value = "20220304-TestAut-PA-35523"
$ws.Range("R13").Value = "'03/04/2022"
$ws.Range("T13").Value = "TestCartridge5523"

